$d = $word.ActiveDocument

# Locate the "Conference Publications" list item that precedes the new
# reference to be inserted (the MSR/PCI 2017 Kechagia reference is the
# last bullet in that sub-list, right before the "Invited Talks" section).
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Analyzing programming languages")) {
        $anchorPara = $p
    }
}

# Insert a brand-new paragraph right after it; InsertParagraphAfter()
# clones the paragraph's formatting (the "Compact" style and the
# numId=3 / ilvl=1 bullet numbering), matching the rest of the list.
$anchorPara.Range.InsertParagraphAfter()

$newPara = $anchorPara.Next()
$newPara.Range.Text = "Stefanos Georgiou, Stamatia Rizou, and Diomidis Spinellis. Software Development Life Cycle for Energy-Efficiency: Techniques and Tools. In ACM Computing Surveys [Submitted for review on May of 2017]"
